$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keeps existing formatting on rows 1-15),
# so the shared-string table gets rebuilt from scratch in the order we
# write values below.
$ws.Range("A1:E15").ClearContents()

# --- Header row ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "from_bus"
$ws.Range("D1").Value = "to_bus"
$ws.Range("E1").Value = "in_service"

# --- Data rows (row, A, B-label, C, D, E) ---
$data = @(
  @(2,  0, "line1", 7,  9,  $true),
  @(3,  1, "line2", 9,  8,  $true),
  @(4,  2, "line3", 8,  10, $true),
  @(5,  3, "line4", 8,  11, $false),
  @(6,  4, "line5", 10, 5,  $true),
  @(7,  5, "line6", 12, 8,  $true),
  @(8,  6, "line7", 14, 11, $true),
  @(9,  7, "line8", 16, 9,  $true),
  @(10, 8, "extr1", 5,  12, $true),
  @(11, 9, "extr2", 5,  9,  $true),
  @(12, 10, "extr3", 10, 11, $false),
  @(13, 11, "extr4", 7,  8,  $true),
  @(14, 12, "extr5", 9,  11, $false),
  @(15, 13, "extr6", 7,  11, $true),
  @(16, 14, "extr7", 5,  7,  $true),
  @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# Rows 16 and 17 are new and fall outside the original formatted range,
# so copy the bold/border/center formatting used by the other "A" column
# cells (and row) onto them.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null
